# Main Table 1 - correct typo ("inclonclusive" -> "inconclusive") and
# store the position value of row 5 (Xq12 / rs189618857:66156010:T:A) as
# text instead of a number, matching the rest of the position column's
# formatting for that SNP string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table1")

# --- Fix the "inclonclusive" -> "inconclusive" typo in the
#     coloc_sexIA column (column U). ---
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 21)
    if ($cell.Text -eq "inclonclusive") {
        $cell.Value = "inconclusive"
    }
}

# --- Store E5's position as text "66156010" instead of a number. ---
$e5 = $ws.Range("E5")
$e5.NumberFormat = "@"
$e5.Value = "66156010"
$e5.Style = "Normal"
